$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44304, 0, 7, 116.4531691898187),
    @(44305, 0, 6, 99.81700216270171),
    @(44306, 2, 5, 83.18083513558476),
    @(44307, 1, 5, 83.18083513558476)
)

$startRow = 230

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy the format from the row above (column A carries the date style)
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = $false
